# Modify command style for both earnings and payment
#
# The sequence diagram used a Laravel-style "name/value" argument notation
# ("paid idx/1 amt/200 m/8 y/2018") in three labels; switch all three to a
# plain positional notation ("paid 1 200 8 2018").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$quote = [char]8220   # “
$rquote = [char]8221  # ”

# --- "TextBox 25": execute( "paid idx/1 amt/200 m/8 y/2018") -------------
$sh = $s.Shapes.Item("TextBox 25")
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Text.Length)
$full.Text = "execute( " + $quote + "paid 1 200 8 2018" + $rquote + ")"

# --- "TextBox 79": parseCommand("paid idx/1 amt/200 m/8 y/2018") ---------
$sh = $s.Shapes.Item("TextBox 79")
$tr = $sh.TextFrame.TextRange
$len = $tr.Text.Length
# Keep the leading "parseCommand" run untouched; rewrite everything after it.
$rest = $tr.Characters(13, $len - 12)
$rest.Text = "(" + $quote + "paid 1 200 8 2018" + $rquote + ")"

# --- "TextBox 60": Parse("idx/1 amt/200 m/8 y/2018") ---------------------
$sh = $s.Shapes.Item("TextBox 60")
$tr = $sh.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Text.Length)
$full.Text = "Parse(" + $quote + "1 200 8 2018" + $rquote + ")"
# The autofit recompute triggered by the text edit above can land a hair
# under the true rendered height due to points<->EMU float rounding;
# nudge it by a sub-point epsilon to settle on the exact rendered value.
$sh.Height = $sh.Height + 0.0001
